$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: replace the text of a Range while forcing Word to keep it as
# a run that is distinct from its neighbours (mirrors what happens when
# a user selects text and retypes over it - a fresh run boundary is
# created even if the resulting formatting is identical to the runs
# around it). We do this by toggling Bold on then off on the freshly
# inserted text.
# ---------------------------------------------------------------------
function Set-RangeTextSplit($rng, $newText) {
    $rng.Text = $newText
    $sub = $d.Range($rng.Start, $rng.Start + $newText.Length)
    $sub.Bold = 1
    $sub.Bold = 0
}

# ---------------------------------------------------------------------
# Helper: merge several adjacent runs that currently hold $fullText
# (spread across more than one <w:r>) into a single run. We overwrite
# the range with a throwaway placeholder (forcing the engine to delete
# the old runs) and then write the real text back into that single
# spot in one shot.
# ---------------------------------------------------------------------
function Merge-RangeRuns($rng) {
    $txt = $rng.Text
    $rng.Text = "PLACEHOLDER_MERGE_TOKEN"
    $sub = $d.Range($rng.Start, $rng.Start + "PLACEHOLDER_MERGE_TOKEN".Length)
    $sub.Text = $txt
}

# 1) "Forum post must be submitted" -> "Forum post " + "is" + " submitted"
$rng = $d.Content
$rng.Find.Execute("Forum post must be submitted", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$s = $rng.Start
$sub = $d.Range($s + 11, $s + 18)   # the "must be" span
Set-RangeTextSplit $sub "is"

# 2) Merge " use case" + ". Actor must have posted previously, so " into one run
$rng = $d.Content
$rng.Find.Execute(" use case. Actor must have posted previously, so ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Merge-RangeRuns $rng

# 3) "Deletion must be confirmed" -> "Forum post is deleted"
$rng = $d.Content
$rng.Find.Execute("Deletion must be confirmed", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = "Forum post is deleted"

# 4) Merge "2.  " + "Student or Lecturer" + " confirm the deletion" into one run
$rng = $d.Content
$rng.Find.Execute("2.  Student or Lecturer confirm the deletion", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Merge-RangeRuns $rng

# 5) Remove the run "Student must have submitted their exam" (leave paragraph empty)
$rng = $d.Content
$rng.Find.Execute("Student must have submitted their exam", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = ""

# 6) Remove the run "Lecturer must confirm the submission" (leave paragraph empty)
$rng = $d.Content
$rng.Find.Execute("Lecturer must confirm the submission", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = ""

# 7) Second occurrence of "Input must be validated" (Update Course use case) -> "Course is updated"
$rng = $d.Content
$rng.Find.Execute("Input must be validated", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2 = $d.Range($rng.End, $d.Content.End)
$rng2.Find.Execute("Input must be validated", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2.Text = "Course is updated"
